# Updated cryptos list on Thu Aug 22 03:36:42 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a text value on a cell while preventing Excel from
# auto-converting numeric-looking strings (e.g. "566.35") into real
# numbers, and without leaving the cell's visible style changed.
function Set-TextValue($Range, $Text) {
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "60.160.95"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.590.11"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "566.35"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "141.60"
$ws.Range("E6").Value = "  -0.64%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.18%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.597"
$ws.Range("E8").Value = "  -0.18%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.609.91"
$ws.Range("E9").Value = "  +0.74%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "6.58"
$ws.Range("E10").Value = "  -0.87%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.12%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +6.24%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -6.59%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.051.42"
$ws.Range("E14").Value = "  +0.42%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Range("D15") "60.174.49"
$ws.Range("E15").Value = "  +1.73%  "

# Row 16 - Avalanche
$ws.Range("E16").Value = "  +1.49%  "

# Row 17 - ShibaInu
Set-TextValue $ws.Range("D17") "0.0000140"
$ws.Range("E17").Value = "  +1.97%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.598.55"
$ws.Range("E18").Value = "  +0.47%  "

# Row 19 - Chainlink
Set-TextValue $ws.Range("D19") "11.34"
$ws.Range("E19").Value = "  +9.39%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +1.83%  "

# Row 21 - BitcoinCash
Set-TextValue $ws.Range("D21") "345.38"
$ws.Range("E21").Value = "  +2.58%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +8.01%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.06%  "

# Row 24 - Polygon
$ws.Range("E24").Value = "  +15.78%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "62.94"
$ws.Range("E25").Value = "  -1.92%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue $ws.Range("D26") "0.995"
$ws.Range("E26").Value = "  -0.03%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  -2.16%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "7.64"
$ws.Range("E28").Value = "  +4.43%  "

# Row 29 - PEPE
Set-TextValue $ws.Range("D29") "0.0₃0780"
$ws.Range("E29").Value = "  +0.84%  "

# Row 30 - PancakeSwap
Set-TextValue $ws.Range("D30") "1.79"
$ws.Range("E30").Value = "  +7.16%  "

# Row 32 - Aptos
$ws.Range("E32").Value = "  +3.58%  "

# Row 33 - Monero
Set-TextValue $ws.Range("D33") "160.80"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "19.43"
$ws.Range("E34").Value = "  +2.64%  "

# Row 35 - NEARProtocol
Set-TextValue $ws.Range("D35") "4.23"
$ws.Range("E35").Value = "  +5.26%  "

# Row 36 - Fetch.AI
Set-TextValue $ws.Range("D36") "0.954"
$ws.Range("E36").Value = "  +9.07%  "

# Row 37 - ImmutableX
Set-TextValue $ws.Range("D37") "1.21"
$ws.Range("E37").Value = "  +4.28%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +7.41%  "

# Row 39 - OKB
Set-TextValue $ws.Range("D39") "37.73"
$ws.Range("E39").Value = "  +0.74%  "

# Row 40 - was SuiNetwork, now Filecoin
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D40") "3.81"
$ws.Range("E40").Value = "  +4.29%  "

# Row 41 - was Filecoin, now SuiNetwork
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D41") "0.852"
$ws.Range("E41").Value = "  -2.64%  "

# Row 42 - Bittensor
Set-TextValue $ws.Range("D42") "293.88"
$ws.Range("E42").Value = "  -0.06%  "

# Row 43 - Aave
Set-TextValue $ws.Range("D43") "138.06"
$ws.Range("E43").Value = "  +4.64%  "

# Row 44 - FirstDigitalUSD
Set-TextValue $ws.Range("D44") "0.997"
$ws.Range("E44").Value = "  -0.21%  "

# Row 45 - was Stellar, now Mantle
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D45") "0.605"
$ws.Range("E45").Value = "  +1.42%  "

# Row 46 - was Mantle, now Stellar
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D46") "0.0977"
$ws.Range("E46").Value = "  +0.45%  "

# Row 47 - EnergySwap
Set-TextValue $ws.Range("D47") "19.54"
$ws.Range("E47").Value = "  +2.53%  "

# Row 48 - Hedera
Set-TextValue $ws.Range("D48") "0.0543"
$ws.Range("E48").Value = "  +1.45%  "

# Row 49 - was VeChain, now RenderToken
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D49") "4.91"
$ws.Range("E49").Value = "  +8.75%  "

# Row 50 - was RenderToken, now VeChain
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0239"
$ws.Range("E50").Value = "  +2.62%  "

# Row 51 - WhiteBITCoin
$ws.Range("E51").Value = "  +0.26%  "
